$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3940.1482
$ws.Range("J40").Value = 3975.8
$ws.Range("L40").Value = 3975.8
$ws.Range("N40").Value = -4325.8
$ws.Range("H43").Value = 3812.3333
$ws.Range("I43").Value = 2612.5
$ws.Range("J43").Value = 4412.25
$ws.Range("K43").Value = 2612.5
$ws.Range("L43").Value = 4412.25
$ws.Range("M43").Value = -2543.5
$ws.Range("N43").Value = -4550.25
$ws.Range("H69").Value = 10013.25
$ws.Range("J69").Value = 10014.968
$ws.Range("L69").Value = 30044.904
$ws.Range("N69").Value = -31792.904
$ws.Range("H72").Value = 10013.25
$ws.Range("J72").Value = 10014.968
$ws.Range("L72").Value = 90134.712
$ws.Range("N72").Value = -98870.712
$ws.Range("H101").Value = 854.125
$ws.Range("J101").Value = 1496.6666
$ws.Range("L101").Value = 4489.9998
$ws.Range("N101").Value = -7733.9998
$ws.Range("H132").Value = 3287.2593
$ws.Range("I132").Value = 3092.0908
$ws.Range("J132").Value = 4146
$ws.Range("K132").Value = 9276.2724
$ws.Range("L132").Value = 12438
$ws.Range("M132").Value = -6746.2724
$ws.Range("N132").Value = -17498
$ws.Range("H138").Value = 2426.21
$ws.Range("I138").Value = 1205.5294
$ws.Range("J138").Value = 2676.229
$ws.Range("K138").Value = 3616.5882
$ws.Range("L138").Value = 8028.687
$ws.Range("M138").Value = 1523.4118
$ws.Range("N138").Value = -18308.687

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4013.4167
$ws.Range("J63").Value = 4832.75
$ws.Range("L63").Value = 4832.75
$ws.Range("N63").Value = -6204.75
$ws.Range("H66").Value = 4013.4167
$ws.Range("J66").Value = 4832.75
$ws.Range("L66").Value = 24163.75
$ws.Range("N66").Value = -31027.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2298.3
$ws.Range("I107").Value = 935
$ws.Range("K107").Value = 935
$ws.Range("M107").Value = 985

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2299.6667
$ws.Range("I2").Value = 2299.6667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2299.6667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2186.6667
$ws.Range("N2").ClearContents()
$ws.Range("H16").Value = 2267.7144
$ws.Range("I16").Value = 2229.1667
$ws.Range("K16").Value = 2229.1667
$ws.Range("M16").Value = -1942.1667
$ws.Range("H42").Value = 2000
$ws.Range("I42").Value = 2000
$ws.Range("K42").Value = 2000
$ws.Range("M42").Value = -1407
$ws.Range("H58").Value = 2880.05
$ws.Range("I58").Value = 2718.3
$ws.Range("J58").Value = 3688.8
$ws.Range("K58").Value = 2718.3
$ws.Range("L58").Value = 3688.8
$ws.Range("M58").Value = -2515.3
$ws.Range("N58").Value = -4094.8
$ws.Range("H59").Value = 47552
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -47290
$ws.Range("H105").Value = 3449.5
$ws.Range("I105").Value = 2765.6667
$ws.Range("J105").Value = 4133.3335
$ws.Range("K105").Value = 2765.6667
$ws.Range("L105").Value = 4133.3335
$ws.Range("M105").Value = -1018.6667
$ws.Range("N105").Value = -7627.3335
$ws.Range("H113").Value = 2267.7144
$ws.Range("I113").Value = 2229.1667
$ws.Range("K113").Value = 2229.1667
$ws.Range("M113").Value = -59.16670000000022
$ws.Range("I122").Value = 3028.2222
$ws.Range("J122").Value = 5927
$ws.Range("K122").Value = 9084.6666
$ws.Range("L122").Value = 17781
$ws.Range("M122").Value = -6634.6666
$ws.Range("N122").Value = -22681
$ws.Range("H132").Value = 4920.4644
$ws.Range("J132").Value = 9105.4
$ws.Range("L132").Value = 27316.2
$ws.Range("N132").Value = -32376.2
$ws.Range("H136").Value = 2880.05
$ws.Range("I136").Value = 2718.3
$ws.Range("J136").Value = 3688.8
$ws.Range("K136").Value = 8154.900000000001
$ws.Range("L136").Value = 11066.4
$ws.Range("M136").Value = -5604.900000000001
$ws.Range("N136").Value = -16166.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H121").Value = 5625770.5
$ws.Range("I121").Value = 517.8333
$ws.Range("J121").Value = 9000922
$ws.Range("K121").Value = 1553.4999
$ws.Range("L121").Value = 27002766
$ws.Range("M121").Value = -243.4999
$ws.Range("N121").Value = -27005386
$ws.Range("H132").Value = 1568.6666
$ws.Range("J132").Value = 1568.6666
$ws.Range("L132").Value = 14117.9994
$ws.Range("N132").Value = -19177.9994
$ws.Range("H137").Value = 8070.222
$ws.Range("I137").Value = 1778.8334
$ws.Range("J137").Value = 11215.917
$ws.Range("K137").Value = 5336.5002
$ws.Range("L137").Value = 33647.751
$ws.Range("M137").Value = -236.5002000000004
$ws.Range("N137").Value = -43847.751
$ws.Range("H140").Value = 2151.6
$ws.Range("I140").Value = 1648.9412
$ws.Range("K140").Value = 4946.8236
$ws.Range("M140").Value = 233.1764000000003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2476.3635
$ws.Range("J126").Value = 2387.0908
$ws.Range("L126").Value = 7161.2724
$ws.Range("N126").Value = -12101.2724
$ws.Range("H134").Value = 87141.42999999999
$ws.Range("J134").Value = 87141.42999999999
$ws.Range("L134").Value = 261424.29
$ws.Range("N134").Value = -266494.29
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 72736.664
$ws.Range("J136").Value = 72736.664
$ws.Range("L136").Value = 218209.992
$ws.Range("N136").Value = -223309.992

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3493.6667
$ws.Range("I22").Value = 1510.25
$ws.Range("J22").Value = 5080.4
$ws.Range("K22").Value = 1510.25
$ws.Range("L22").Value = 5080.4
$ws.Range("M22").Value = -1215.25
$ws.Range("N22").Value = -5670.4
$ws.Range("H27").Value = 3493.6667
$ws.Range("I27").Value = 1510.25
$ws.Range("J27").Value = 5080.4
$ws.Range("K27").Value = 1510.25
$ws.Range("L27").Value = 5080.4
$ws.Range("M27").Value = -1403.25
$ws.Range("N27").Value = -5294.4
$ws.Range("H46").Value = 7246
$ws.Range("I46").Value = 2566.6667
$ws.Range("J46").Value = 8415.833000000001
$ws.Range("K46").Value = 2566.6667
$ws.Range("L46").Value = 8415.833000000001
$ws.Range("M46").Value = -2378.6667
$ws.Range("N46").Value = -8791.833000000001
$ws.Range("H55").Value = 337.42856
$ws.Range("I55").Value = 240.25
$ws.Range("K55").Value = 240.25
$ws.Range("M55").Value = -67.25
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 5666.143
$ws.Range("I132").Value = 5390.4
$ws.Range("J132").Value = 5916.8184
$ws.Range("K132").Value = 16171.2
$ws.Range("L132").Value = 17750.4552
$ws.Range("M132").Value = -13641.2
$ws.Range("N132").Value = -22810.4552
$ws.Range("H136").Value = 6874.857
$ws.Range("I136").Value = 4579.5625
$ws.Range("J136").Value = 14219.8
$ws.Range("K136").Value = 13738.6875
$ws.Range("L136").Value = 42659.39999999999
$ws.Range("M136").Value = -11188.6875
$ws.Range("N136").Value = -47759.39999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5064.2856
$ws.Range("J62").Value = 6425
$ws.Range("L62").Value = 6425
$ws.Range("N62").Value = -7673
$ws.Range("H65").Value = 5064.2856
$ws.Range("J65").Value = 6425
$ws.Range("L65").Value = 32125
$ws.Range("N65").Value = -38365
$ws.Range("H96").Value = 12326.733
$ws.Range("I96").Value = 6444.222
$ws.Range("K96").Value = 6444.222
$ws.Range("M96").Value = -5071.222
$ws.Range("H100").Value = 959.13336
$ws.Range("I100").Value = 505.375
$ws.Range("J100").Value = 1477.7142
$ws.Range("K100").Value = 1010.75
$ws.Range("L100").Value = 2955.4284
$ws.Range("M100").Value = -469.75
$ws.Range("N100").Value = -4037.4284
$ws.Range("H122").Value = 5678.6523
$ws.Range("J122").Value = 8814.25
$ws.Range("L122").Value = 26442.75
$ws.Range("N122").Value = -31342.75
$ws.Range("H137").Value = 52500
$ws.Range("J137").Value = 52500
$ws.Range("L137").Value = 52500
$ws.Range("N137").Value = -62700
